# UniformA-HW25: "Updated notebook, reran simulation"
#
# The refreshed simulation grew the "Kind" category list (column B) from 28
# to 30 entries:
#   - two new entries, "Holden" and "Rizzie Spiral", were inserted right
#     after "Spiral5" (pushing every later entry down by two rows), and
#   - "Thomas Hex" was renamed to "Matthies Hex".
#   - the two entries that fell off the end of the old range ("Michael-CCHex",
#     "Michael-SNHex") re-appear as two brand-new trailing rows (30 and 31),
#     each filled with the same per-column tally of 1 as every other row.
#
# Column headers (row 1 and row 2) and the per-cell tallies (columns C:W)
# are unchanged in meaning, so only column B needs to be rewritten for the
# rows that shifted; rows 30-31 are appended in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing category labels in column B down by two rows ---
# (rows 4-29 each take on the label that used to belong two rows above them;
# row 11 picks up the renamed "Matthies Hex" label in the process)
$ws.Range("B4").Value = "Holden"
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("B6").Value = "RotRing OmegaMax-90"
$ws.Range("B7").Value = "Equal Angle"
$ws.Range("B8").Value = "Tilt Rotate"
$ws.Range("B9").Value = "CLR"
$ws.Range("B10").Value = "Rizzie Hex"
$ws.Range("B11").Value = "Matthies Hex"
$ws.Range("B12").Value = "Tilt Rotate_Partial"
$ws.Range("B13").Value = "RotRing OmegaMax-60"
$ws.Range("B14").Value = "Equal Angle_Partial"
$ws.Range("B15").Value = "Rizzie Hex_Partial"
$ws.Range("B16").Value = "ND Single"
$ws.Range("B17").Value = "RD Single"
$ws.Range("B18").Value = "TD Single"
$ws.Range("B19").Value = "Morris Single"
$ws.Range("B20").Value = "Ring Perpendicular to ND"
$ws.Range("B21").Value = "Ring Perpendicular to RD"
$ws.Range("B22").Value = "Ring Perpendicular to TD"
$ws.Range("B23").Value = "OffsetFTD"
$ws.Range("B24").Value = "OffsetATD"
$ws.Range("B25").Value = "OffsetF45"
$ws.Range("B26").Value = "OffsetA45"
$ws.Range("B27").Value = "OffsetFRD"
$ws.Range("B28").Value = "OffsetARD"
$ws.Range("B29").Value = "Gaussian Quadrature"

# --- Append the two new simulation rows (30 and 31) ---
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("C30:W30").Value = 1

$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Michael-SNHex"
$ws.Range("C31:W31").Value = 1
